$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Issue #13: allow two columns in metadata files to be related to create
# hierarchical SKOS. A new row (row 2) of "slug" identifiers is inserted
# right below the human-readable header row, giving every header column a
# machine-friendly key that other columns can reference.
$ws.Rows.Item(2).Insert()

$ws.Range("A2").Value = "n-parados"
$ws.Range("B2").Value = "comarca-nombre"
$ws.Range("C2").Value = "comarca-codigo"
$ws.Range("D2").Value = "aragon"
$ws.Range("E2").Value = "provincia-codigo"
$ws.Range("F2").Value = "provincia-nombre"
$ws.Range("G2").Value = "sexo"
$ws.Range("H2").Value = "mes-y-ano"
